# Convert the old "bookmarked Heading1 title + bold 'By ' byline" opening
# into a pandoc-style title block: a Title-styled paragraph with the bare
# title (run-per-word, matching pandoc's docx writer) followed by an
# Authors-styled paragraph with just the author name (no "By " prefix).

$d = $word.ActiveDocument

# The document opens with:
#   1) Heading1 paragraph:  "Day After Day - June 1943"
#   2) Normal/bold paragraph: "By Dorothy Day"
# Grab both as a single Range so the replacement lands in one shot.
$p1 = $d.Paragraphs(1)
$p2 = $d.Paragraphs(2)

if ($p1.Range.Text.TrimEnd([char]13) -ne "Day After Day - June 1943") {
    throw "Unexpected paragraph 1 text: $($p1.Range.Text)"
}
if ($p2.Range.Text.TrimEnd([char]13) -ne "By Dorothy Day") {
    throw "Unexpected paragraph 2 text: $($p2.Range.Text)"
}

$titleRange = $d.Range($p1.Range.Start, $p2.Range.End)

$openXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr><w:pStyle w:val="Title"/></w:pPr>
            <w:r><w:t xml:space="preserve">Day</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:t xml:space="preserve">After</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:t xml:space="preserve">Day</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:t xml:space="preserve">-</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:t xml:space="preserve">June</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:t xml:space="preserve">1943</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr><w:pStyle w:val="Authors"/></w:pPr>
            <w:r><w:t xml:space="preserve">Dorothy</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:t xml:space="preserve">Day</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$titleRange.InsertXML($openXml)
